$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dist_coûts")

# Header row: Z1 and AA1 become F2 / F3 (Y1 stays F1)
$ws.Range("Z1").Value = "F2"
$ws.Range("AA1").Value = "F3"

# Rows 2-6: update Y/Z/AA values
$ws.Range("Y2").Value = 2.25
$ws.Range("Z2").Value = 50
$ws.Range("AA2").Value = 0

$ws.Range("Y3").Value = 2.25
$ws.Range("Z3").Value = 50
$ws.Range("AA3").Value = 0

$ws.Range("Y4").Value = 2.25
$ws.Range("Z4").Value = 55
$ws.Range("AA4").Value = 0

$ws.Range("Y5").Value = 2.25
$ws.Range("Z5").Value = 60
$ws.Range("AA5").Value = 0

$ws.Range("Y6").Value = 2.25
$ws.Range("Z6").Value = 60
$ws.Range("AA6").Value = 0

# Row 11: new values
$ws.Range("W11").Value = "Spécifique"
$ws.Range("Y11").Value = 1
$ws.Range("Z11").Value = 0
$ws.Range("AA11").Value = 0

$wb.Save()
